$wb = $excel.ActiveWorkbook
$north = $wb.Worksheets.Item(1)
$south = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# North sheet ("North" / North_cons table) - revised "Dates" column text
# ---------------------------------------------------------------------------
$north.Range("B2").Value = "May 14 - May 23, June 13 - June 22, July 12 - July 21"
$north.Range("B3").Value = "August 10 - August 19, September 9 -September 18, October 8 - October 17"
$north.Range("B4").Value = "February 14 - February 23, March 14 - March 24"
$north.Range("B5").Value = "June 13 - June 22, July 12 - July 21, August 10 - August 19"
$north.Range("B6").Value = "April 14-23, May 14-23"
$north.Range("B7").Value = "January 16 - January 25, February 14 - February 23, March 14 - March 24"
$north.Range("B8").Value = "October 8 - October17, November 7 - November 16"
$north.Range("B9").Value = "January 16 - January 25, November 7- November 16, December 6 - December 15"
$north.Range("B10").Value = "January 16 - January 25"

# Rows whose wrapped text now spans two lines get a taller row height
$north.Rows.Item(3).RowHeight = 26.4
$north.Rows.Item(5).RowHeight = 26.4
$north.Rows.Item(7).RowHeight = 26.4
$north.Rows.Item(9).RowHeight = 26.4

# ---------------------------------------------------------------------------
# South sheet ("South" / South_cons table) - revised "Dates" column text
# ---------------------------------------------------------------------------
$south.Range("B2").Value = "June 13- June 22"
$south.Range("B3").Value = "February 14- February 23, March 14 - March 24"
$south.Range("B4").Value = "April 14 - April 23,  May 14 - May 23, June 13 - June 22"
$south.Range("B5").Value = "September 9- September 18, October 8- October 17, November 7- November 16, December 6 -December 15"
$south.Range("B6").Value = "July 12 - July 21"
$south.Range("B7").Value = "April 14- April 23, May 14 - May 23"
$south.Range("B8").Value = "January 16- January 25, February 14 -February 23, March 14 - March 24"
$south.Range("B9").Value = "October 8 - October 17, November 7 - November 16"
$south.Range("B10").Value = "August 10 - August 19, September 9 -September 18"
$south.Range("B11").Value = "July 12 - July 21, August 10 - August 19"
$south.Range("B12").Value = "January 16 - January 25"

# Rows whose wrapped text now spans two/three lines get a taller row height
$south.Rows.Item(4).RowHeight = 26.4
$south.Rows.Item(5).RowHeight = 39.6
$south.Rows.Item(8).RowHeight = 26.4

# ---------------------------------------------------------------------------
# Restore selections left behind by the editor
# ---------------------------------------------------------------------------
$north.Activate() | Out-Null
$north.Range("B13").Select() | Out-Null

$south.Activate() | Out-Null
$south.Range("F8").Select() | Out-Null
